$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2:O35").Value = "2022-08-27 20:57:52"
